$d = $word.ActiveDocument

# Locate the block that starts at the "Table " paragraph and runs through
# the page-break paragraph right before "Description des données sur les
# tables à créer". We find it by searching for the literal text so this
# does not depend on brittle character offsets.

$startFind = $d.Content.Duplicate
$startFind.Find.Execute("Table ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPara = $startFind.Paragraphs(1)

$endFind = $d.Content.Duplicate
$endFind.Find.Execute("Description des données sur les tables à créer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$blockStart = $startPara.Range.Start
$blockEnd = $endFind.Paragraphs(1).Previous().Range.End

$target = $d.Range($blockStart, $blockEnd)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
'<w:p><w:pPr><w:spacing w:after="120"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr></w:p>' + `
'<w:p><w:pPr><w:spacing w:after="120"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr></w:p>' + `
'<w:p><w:pPr><w:spacing w:after="120"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr></w:p>' + `
'<w:p><w:pPr><w:spacing w:after="120"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr></w:p>' + `
'<w:p><w:pPr><w:spacing w:after="120"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr></w:p>' + `
'<w:p><w:pPr><w:spacing w:after="120"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr></w:p>' + `
'<w:p><w:pPr><w:spacing w:after="120"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr></w:p>' + `
'<w:p><w:pPr><w:spacing w:after="160" w:line="259" w:lineRule="auto"/><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="24"/></w:rPr><w:br w:type="page"/></w:r></w:p>' + `
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)

Write-Output "done"
